# "Actualizacion base de datos" - refresh the hitos (milestones) table:
#   - header: month | year | type | link   (was: month | year | milestone | type)
#   - replace all data rows with the refreshed dataset (20 rows)
#   - the single DOI-style link that still reads as a real hyperlink (D5) gets
#     an actual Hyperlink + the built-in Hyperlink style
#   - column B/C get explicit widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe the old table (A1:D6) so no stale cells survive -------------------
$ws.Range("A1:D6").Clear()

# --- header row --------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "month"
$ws.Cells.Item(1,2).Value = "year"
$ws.Cells.Item(1,3).Value = "type"
$ws.Cells.Item(1,4).Value = "link"

# --- data rows ----------------------------------------------------------------
$rows = @(
  @(1,  2017, "Creación del grupo de Ecoinformática", ""),
  @(2,  2022, "Seminarios", ""),
  @(10, 2024, "Primeras Jornadas Ecoinformáticas ", "https://ecoinf.quarto.pub/iecoinf/"),
  @(3,  2017, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2017.26-1.20"),
  @(6,  2017, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2017.26-2.08"),
  @(11, 2017, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2017.26-3.14"),
  @(3,  2018, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1570"),
  @(7,  2018, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1604"),
  @(6,  2018, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1604"),
  @(3,  2019, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1699"),
  @(11, 2019, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1880"),
  @(2,  2020, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1948"),
  @(10, 2019, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1838"),
  @(5,  2020, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.1995"),
  @(12, 2020, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2129"),
  @(6,  2021, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2256"),
  @(12, 2021, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2332"),
  @(2,  2023, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2527"),
  @(6,  2024, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2797"),
  @(8,  2024, "Nota ecoinformática", "https://doi.org/10.7818/ECOS.2745")
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  if ($row[3] -ne "") {
    $ws.Cells.Item($r,4).Value = $row[3]
  }
  $r = $r + 1
}

# --- the one cell that is a real hyperlink (D5) -------------------------------
$ws.Hyperlinks.Add($ws.Range("D5"), "https://doi.org/10.7818/ECOS.2017.26-1.20") | Out-Null

# --- column widths -------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19.7265625
$ws.Columns.Item(3).ColumnWidth = 38.36328125

# --- final selection (matches the author's last click before saving) --------
$ws.Range("J8").Select() | Out-Null
